# Update the "Metadata" sheet values for Title, Date, and Description.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Title value (row 5, column B)
$ws.Range("B5").Value = "Modèle logique métier - FR LM Non remboursable"

# Date value (row 8, column B)
$ws.Range("B8").Value = "2026-01-16T13:49:34+00:00"

# Description value (row 12, column B)
$ws.Range("B12").Value = "Non remboursable"
